$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item("login")

# Add the new "signUp" worksheet right after "login"
$newSheet = $wb.Worksheets.Add($null, $loginSheet)
$newSheet.Name = "signUp"

# ---- Header row (bold, filled) ----
$newSheet.Range("A1").Value = "test step"
$newSheet.Range("B1").Value = "test step"
$newSheet.Range("C1").Value = "test step"
$newSheet.Range("D1").Value = "test step"
$newSheet.Range("E1").Value = "test step"

# ---- Row 2 ----
$newSheet.Range("A2").Value = "open browser"
$newSheet.Range("B2").Value = "NA"
$newSheet.Range("C2").Value = "NA"
$newSheet.Range("D2").Value = "open browser"
$newSheet.Range("E2").Value = "chrome"

# ---- Row 3 ----
$newSheet.Range("A3").Value = "launch url"
$newSheet.Range("B3").Value = "NA"
$newSheet.Range("C3").Value = "NA"
$newSheet.Range("D3").Value = "enter url"
$newSheet.Range("E3").Value = "http://app.hubspot.com/login"

# ---- Row 4 ----
$newSheet.Range("A4").Value = "verify sign up link"
$newSheet.Range("B4").Value = "linkText"
$newSheet.Range("C4").Value = "Sign up"
$newSheet.Range("D4").Value = "click"
$newSheet.Range("E4").Value = "NA"

# ---- Row 5 ----
$newSheet.Range("A5").Value = "close browser"
$newSheet.Range("B5").Value = "NA"
$newSheet.Range("C5").Value = "NA"
$newSheet.Range("D5").Value = "quit"
$newSheet.Range("E5").Value = "NA"

# ---- Hyperlink on E3 ----
$e3 = $newSheet.Range("E3")
$newSheet.Hyperlinks.Add($e3, "http://app.hubspot.com/login", "", "", "http://app.hubspot.com/login") | Out-Null

# ---- Fonts: body text is Arial 12 black ----
$body = $newSheet.Range("A1:E5")
$body.Font.Name = "Arial"
$body.Font.Size = 12
$body.Font.Color = 0

# ---- Header row: bold + fill ----
$header = $newSheet.Range("A1:E1")
$header.Font.Bold = $true
$header.Interior.Color = 2381624

# ---- Hyperlink cell formatting (underline, Excel hyperlink blue) ----
$e3.Font.Underline = 2
$e3.Font.Color = 12673797

# ---- Row heights ----
$newSheet.Rows.Item(1).RowHeight = 15.75
$newSheet.Rows.Item(2).RowHeight = 15.75
$newSheet.Rows.Item(3).RowHeight = 15.75
$newSheet.Rows.Item(4).RowHeight = 15.75
$newSheet.Rows.Item(5).RowHeight = 15.75

# ---- Column widths (approximate target character widths) ----
$newSheet.Columns.Item(1).ColumnWidth = 19.15625
$newSheet.Columns.Item(2).ColumnWidth = 20.91796875
$newSheet.Columns.Item(3).ColumnWidth = 19.53515625
$newSheet.Columns.Item(4).ColumnWidth = 20.7734375
$newSheet.Columns.Item(5).ColumnWidth = 33.53515625

# ---- Selection on the new sheet, then return focus to "login" ----
$newSheet.Range("B19").Select() | Out-Null
$loginSheet.Select() | Out-Null
